# modifiche float e outstanding da finviz
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (DWTX) updates ---
$ws.Range("H2").Value = 8.25
$ws.Range("O2").Value = 49.46
$ws.Range("Q2").Value = -7.25
$ws.Range("AJ2").Value = 8.33
$ws.Range("AK2").Value = 8.220000000000001
$ws.Range("AN2").Value = 8.33
$ws.Range("AO2").Value = 8.220000000000001

# --- Row 3 (JFB) updates ---
$ws.Range("D3").Value = 3010000
$ws.Range("H3").Value = 10.72
$ws.Range("O3").Value = 56.5
$ws.Range("P3").Value = 11.55
$ws.Range("Q3").Value = 7.79
$ws.Range("AJ3").Value = 10.9
$ws.Range("AK3").Value = 10.69
$ws.Range("AN3").Value = 10.9
$ws.Range("AO3").Value = 10.69

# --- Row 4 (MSS) updates ---
$ws.Range("D4").Value = 2930000
$ws.Range("H4").Value = 3.75
$ws.Range("O4").Value = 292.67
$ws.Range("Q4").Value = -56.31
$ws.Range("T4").Value = 3.93
$ws.Range("X4").Value = 3.93
$ws.Range("AB4").Value = 3.93
$ws.Range("AF4").Value = 3.93
$ws.Range("AJ4").Value = 3.75
$ws.Range("AN4").Value = 3.75

# --- Row 5 (POAI) updates ---
$ws.Range("H5").Value = 2
$ws.Range("O5").Value = 158.06
$ws.Range("Q5").Value = -25.56
$ws.Range("T5").Value = 2.03
$ws.Range("X5").Value = 2.03
$ws.Range("AB5").Value = 2.03
$ws.Range("AF5").Value = 2.03
$ws.Range("AJ5").Value = 2
$ws.Range("AK5").Value = 1.92
$ws.Range("AN5").Value = 2
$ws.Range("AO5").Value = 1.92

# --- Insert a new row 6 (shifts old row 6 ZURA down to row 7) ---
$ws.Rows.Item(6).Insert()

# --- Row 6: new ticker RDHL ---
$ws.Range("A6").Value = "RDHL"
$ws.Range("B6").Value = 7791872
$ws.Range("C6").Value = 3329860
$ws.Range("D6").Value = 3310000
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.0269
$ws.Range("G6").Value = 45929
$ws.Range("H6").Value = 2.5
$ws.Range("I6").Value = 1.84
$ws.Range("J6").Value = 2.79
$ws.Range("K6").Value = 2.01
$ws.Range("L6").Value = 2.34
$ws.Range("M6").Value = 55704600
$ws.Range("N6").Value = 2.58
$ws.Range("O6").Value = 35.87
$ws.Range("P6").Value = 2.34
$ws.Range("Q6").Value = -6.29
$ws.Range("R6").Value = ""
$ws.Range("S6").Value = ""
$ws.Range("T6").Value = 2.79
$ws.Range("U6").Value = 2.39
$ws.Range("V6").Value = 22893516
$ws.Range("W6").Value = "no"
$ws.Range("X6").Value = 2.79
$ws.Range("Y6").Value = 2.18
$ws.Range("Z6").Value = 50923227
$ws.Range("AA6").Value = "no"
$ws.Range("AB6").Value = 2.79
$ws.Range("AC6").Value = 2.06
$ws.Range("AD6").Value = 51481031
$ws.Range("AE6").Value = "no"
$ws.Range("AF6").Value = 2.79
$ws.Range("AG6").Value = 2.04
$ws.Range("AH6").Value = 79789126
$ws.Range("AI6").Value = "no"
$ws.Range("AJ6").Value = 2.54
$ws.Range("AK6").Value = 2.5
$ws.Range("AL6").Value = 17197546
$ws.Range("AM6").Value = "n/a"
$ws.Range("AN6").Value = 2.54
$ws.Range("AO6").Value = 2.5
$ws.Range("AP6").Value = 17197546
$ws.Range("AQ6").Value = "n/a"

# --- Row 7: ZURA (shifted down), with a few changed values ---
$ws.Range("H7").Value = 3.5
$ws.Range("O7").Value = 45.23
$ws.Range("Q7").Value = -7.35
